# Slide Deck.pptx edit — "Add Timestamp and start parsing Victron Data"
#
# 1. Slide 2 ("Content Placeholder 2" outline text): rename the
#    "Daily Energy Data" bullet to "Energy Meter Data".
# 2. Slide 2: add two small text boxes ("+" and "=") that sit between
#    the two screenshots, turning them into a little "A + B = C" diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. Fix up the outline text -----------------------------------------
# Only the "Daily Energy Data" substring changes; everything else in the
# paragraph (tabs, the trailing empty paragraph, run formatting) must stay
# untouched, so we patch just the affected characters instead of
# reassigning the whole TextRange.
$outline = $s.Shapes.Item(2)
$tr = $outline.TextFrame.TextRange
$fullText = $tr.Text
$oldPhrase = "Daily Energy Data"
$newPhrase = "Energy Meter Data"
$startPos = $fullText.IndexOf($oldPhrase)
if ($startPos -ge 0) {
    $tr.Characters($startPos + 1, $oldPhrase.Length).Text = $newPhrase
}

# --- 2. Add the "+" / "=" connector text boxes --------------------------
function Add-ConnectorBox([double]$left, [double]$top, [double]$width, [double]$height, [string]$text) {
    $tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = "Content Placeholder 2"

    $tb.TextFrame.TextRange.Text = $text

    # Match the body-text insets / autofit used for the rest of the deck.
    $tb.TextFrame.MarginLeft = 7.2
    $tb.TextFrame.MarginTop = 3.6
    $tb.TextFrame.MarginRight = 7.2
    $tb.TextFrame.MarginBottom = 3.6
    $tb.TextFrame.Orientation = 1
    $tb.TextFrame.AutoSize = 2

    # These are plain text boxes, so turn off the inherited bullet.
    $tb.TextFrame.TextRange.ParagraphFormat.Bullet.Font.Name = "Arial"
    $tb.TextFrame.TextRange.ParagraphFormat.Bullet.Visible = 0

    return $tb
}

# Sits just to the right of "Picture 7" (the first chart image).
$plusBox = Add-ConnectorBox 325.7476377952756 294.55102362204724 34.25244094488189 28.68755905511811 "+"

# Sits just to the right of "Picture 2" (the second chart image).
$equalsBox = Add-ConnectorBox 611.0268503937008 294.55102362204724 34.25244094488189 28.68755905511811 "="
